$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44511
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("S2").Value = 2800

# Row 5
$ws.Range("D5").Value = 44476
$ws.Range("M5").Value = 120

# Row 6
$ws.Range("D6").Value = 44466
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44473
$ws.Range("M7").Value = 180

# Row 8
$ws.Range("D8").Value = 44434
$ws.Range("M8").Value = 20

# Row 9
$ws.Range("D9").Value = 44435
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 2000

# Row 10 (changed)
$ws.Range("D10").Value = 44517
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 27000
$ws.Range("O10").Value = 27000
$ws.Range("P10").Value = 27000
$ws.Range("S10").Value = 2700

# Row 11 (new)
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = 44517
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = "Tropicales y subtropicales"
$ws.Range("I11").Value = 100108004
$ws.Range("J11").Value = "Papaya"
$ws.Range("K11").Value = "Cultivar IV Región"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("Q11").Value = "$/bandeja 10 kilos"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 2500
$ws.Range("T11").Value = 10

# Row 12 (new)
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44432
$ws.Range("D12").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100108
$ws.Range("H12").Value = "Tropicales y subtropicales"
$ws.Range("I12").Value = 100108004
$ws.Range("J12").Value = "Papaya"
$ws.Range("K12").Value = "Cultivar IV Región"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 2000
$ws.Range("T12").Value = 10
